$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look numeric (single "." decimal) need to be
# forced to Text format first, otherwise Excel auto-converts them to
# numbers (losing the original textual formatting/precision).
$textCells = @(
    "D4",
    "D5",
    "D7",
    "D10",
    "D11",
    "D13",
    "D14",
    "D15",
    "D17",
    "D18",
    "D19",
    "D21",
    "D24",
    "D25",
    "D27",
    "D28",
    "D29",
    "D30",
    "D31",
    "D32",
    "D34",
    "D35",
    "D36",
    "D37",
    "D38",
    "D39",
    "D40",
    "D41",
    "D42",
    "D43",
    "D44",
    "D47",
    "D48",
    "D49",
    "D51",
)
foreach ($c in $textCells) {
    $ws.Range($c).NumberFormat = "@"
}

# Apply the updated price values (text-protected cells)
$ws.Range("D4").Value = "1.002"
$ws.Range("D5").Value = "307.56"
$ws.Range("D7").Value = "0.5064"
$ws.Range("D10").Value = "0.8867"
$ws.Range("D11").Value = "20.63"
$ws.Range("D13").Value = "0.07558"
$ws.Range("D14").Value = "5.317"
$ws.Range("D15").Value = "89.17"
$ws.Range("D17").Value = "0.000008471"
$ws.Range("D18").Value = "14.11"
$ws.Range("D19").Value = "1.001"
$ws.Range("D21").Value = "5.072"
$ws.Range("D24").Value = "6.481"
$ws.Range("D25").Value = "150.58"
$ws.Range("D27").Value = "17.95"
$ws.Range("D28").Value = "2.092"
$ws.Range("D29").Value = "112.69"
$ws.Range("D30").Value = "4.750"
$ws.Range("D31").Value = "4.691"
$ws.Range("D32").Value = "0.09048"
$ws.Range("D34").Value = "3.095"
$ws.Range("D35").Value = "1.159"
$ws.Range("D36").Value = "0.7386"
$ws.Range("D37").Value = "0.02034"
$ws.Range("D38").Value = "2.496"
$ws.Range("D39").Value = "3.044"
$ws.Range("D40").Value = "1.079"
$ws.Range("D41").Value = "0.5346"
$ws.Range("D42").Value = "6.594"
$ws.Range("D43").Value = "115.66"
$ws.Range("D44").Value = "8.317"
$ws.Range("D47").Value = "0.4632"
$ws.Range("D48").Value = "9.944"
$ws.Range("D49").Value = "1.565"
$ws.Range("D51").Value = "36.44"

# Restore default (Normal) style now that the text is safely stored
foreach ($c in $textCells) {
    $ws.Range($c).Style = "Normal"
}

# Remaining plain value updates (coin names, links, prices with
# thousand-separator dots, and volume percentages)
$ws.Range("D2").Value = "27.103.94"
$ws.Range("E2").Value = "  -3.01%  "
$ws.Range("D3").Value = "1.868.13"
$ws.Range("E3").Value = "  -2.14%  "
$ws.Range("E4").Value = "  +0.23%  "
$ws.Range("E5").Value = "  -1.92%  "
$ws.Range("E6").Value = "  +0.22%  "
$ws.Range("E7").Value = "  +0.97%  "
$ws.Range("E8").Value = "  -2.09%  "
$ws.Range("E9").Value = "  -2.15%  "
$ws.Range("E10").Value = "  -2.80%  "
$ws.Range("E11").Value = "  -2.78%  "
$ws.Range("D12").Value = "1.863.55"
$ws.Range("E12").Value = "  -2.44%  "
$ws.Range("E13").Value = "  -1.50%  "
$ws.Range("E14").Value = "  -3.54%  "
$ws.Range("E15").Value = "  -3.84%  "
$ws.Range("E17").Value = "  -3.15%  "
$ws.Range("E18").Value = "  -3.99%  "
$ws.Range("E19").Value = "  +0.27%  "
$ws.Range("D20").Value = "27.158.02"
$ws.Range("E20").Value = "  -2.96%  "
$ws.Range("E21").Value = "  -2.20%  "
$ws.Range("D22").Value = "2.092.52"
$ws.Range("E22").Value = "  -2.28%  "
$ws.Range("E23").Value = "  -2.77%  "
$ws.Range("E24").Value = "  -1.92%  "
$ws.Range("E25").Value = "  -1.57%  "
$ws.Range("E26").Value = "  -0.46%  "
$ws.Range("E27").Value = "  -2.60%  "
$ws.Range("E28").Value = "  -5.35%  "
$ws.Range("E29").Value = "  -2.47%  "
$ws.Range("E30").Value = "  -3.73%  "
$ws.Range("E31").Value = "  -3.47%  "
$ws.Range("E32").Value = "  +0.26%  "
$ws.Range("E33").Value = "  -3.15%  "
$ws.Range("E34").Value = "  -3.61%  "
$ws.Range("B35").Value = "ARBITRUM"
$ws.Range("C35").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("E35").Value = "  -6.24%  "
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("E36").Value = "  -5.20%  "
$ws.Range("E37").Value = "  -2.66%  "
$ws.Range("E38").Value = "  -3.57%  "
$ws.Range("E40").Value = "  -1.44%  "
$ws.Range("E41").Value = "  -3.78%  "
$ws.Range("E42").Value = "  -4.39%  "
$ws.Range("E43").Value = "  +1.95%  "
$ws.Range("E44").Value = "  -2.46%  "
$ws.Range("E45").Value = "  -3.14%  "
$ws.Range("E46").Value = "  +0.23%  "
$ws.Range("E47").Value = "  -4.24%  "
$ws.Range("E48").Value = "  -6.60%  "
$ws.Range("E49").Value = "  -4.56%  "
$ws.Range("E51").Value = "  -1.86%  "

